$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert a brand-new row 7 ("Programowanie") right above the former row 7
#    ("Przelutowac pajeczaka..."), pushing everything below down by one.
#    Excel auto-extends the merged ranges B3:B13 -> B3:B14 and
#    C8:C12 -> C9:C13 because the insertion point sits inside those merges.
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).Insert()

# Former row 6 ("Podlaczyc scp + raspberry przez Serial") becomes the new ESP
# connection task; give it the alternating highlighted look used by sibling
# rows (copy full format from C3, which already carries that style).
$ws.Range("C3").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("C6").Value = "Podłączyć esp + raspberry przez Serial"

# New row 7 content: "Programowanie" header (left/center, wrapped, like the
# plain text cells) + long ESP/Raspberry note (centered, wrapped, like the
# other free-form note cells e.g. D17/D18).
$ws.Range("C4").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = "Programowanie"

$ws.Range("D17").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = "ESP: cleanup kodu, użycie hala i dwóch silników`nRaspberry: uzgodnić układ współrzędnych (środek, orientacja); refactor pathGenerator; klasa do obsługi arm-gimbala"

$ws.Rows.Item(7).RowHeight = 75

# F4: a blank placeholder cell carrying a date number-format, matching the
# month/year tracking column introduced alongside the new sections.
$ws.Range("F4").NumberFormat = "mmm-yy"

# ---------------------------------------------------------------------------
# 2) Append two brand-new rows at the bottom of the list for the new
#    "Naprawy" (Repairs) section, and extend the B23:B24 merge to B23:B26.
# ---------------------------------------------------------------------------
$ws.Range("C25").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = "Zrobić dystans na silnik poziomy"

$ws.Range("C25").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = "Naprawić/przeprojektować mechanizm celowania laserem."

$ws.Range("B27").HorizontalAlignment = -4108
$ws.Range("B27").VerticalAlignment = -4108
$ws.Range("B27").WrapText = $true
$ws.Range("B27").Value = "Naprawy"
$ws.Rows.Item(27).RowHeight = 30

$ws.Range("B23:B24").UnMerge()
$ws.Range("B23:B26").Merge()
$ws.Range("B23:B26").HorizontalAlignment = -4108
$ws.Range("B23:B26").VerticalAlignment = -4108
$ws.Range("B23:B26").WrapText = $true

# ---------------------------------------------------------------------------
# 3) Final view state, matching the author's last on-screen selection.
# ---------------------------------------------------------------------------
$ws.Range("L7").Select()
